$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expedia")

# Write row 2 (data row) first so the shared-string table picks up
# "NewYork", "October", "November" before the header labels, matching
# the order Excel produced when the row was filled in.
$ws.Range("A2").Value = "NewYork"
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = "October"
$ws.Range("D2").Value = 2021
$ws.Range("E2").Value = 21
$ws.Range("F2").Value = "November"
$ws.Range("G2").Value = 2021

# Then the header row.
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "Checkin Day"
$ws.Range("C1").Value = "Checkin Month"
$ws.Range("D1").Value = "Checkin Year"
$ws.Range("E1").Value = "Checkout Day"
$ws.Range("F1").Value = "Checkout Month"
$ws.Range("G1").Value = "CheckoutYear"
